$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = "2026-02-05 22:17:46"
$ws.Range("E3").Value2 = "2026-02-05 22:17:48"
$ws.Range("O3").Value2 = "-1.7 °C"
$ws.Range("E4").Value2 = "2026-02-05 22:17:51"
$ws.Range("O4").Value2 = "11.8 °C"
$ws.Range("E5").Value2 = "2026-02-05 22:17:53"
$ws.Range("O5").Value2 = "10.2 °C"
$ws.Range("E6").Value2 = "2026-02-05 22:17:56"
$ws.Range("H6").Formula = "'69%"
$ws.Range("O6").Value2 = "13.2 °C"
$ws.Range("E7").Value2 = "2026-02-05 22:17:59"
$ws.Range("O7").Value2 = "10.5 °C"
$ws.Range("E8").Value2 = "2026-02-05 22:18:01"
$ws.Range("H8").Formula = "'83%"
$ws.Range("K8").Value2 = "5.5 MJ/m2"
$ws.Range("O8").Value2 = "9.2 °C"
$ws.Range("E9").Value2 = "2026-02-05 22:18:04"
$ws.Range("E10").Value2 = "2026-02-05 22:18:06"
$ws.Range("O10").Value2 = "7.8 °C"
$ws.Range("E11").Value2 = "2026-02-05 22:18:09"
$ws.Range("H11").Formula = "'93%"
$ws.Range("O11").Value2 = "1.0 °C"
$ws.Range("E12").Value2 = "2026-02-05 22:18:11"
$ws.Range("O12").Value2 = "10.6 °C"
$ws.Range("E13").Value2 = "2026-02-05 22:18:14"
$ws.Range("O13").Value2 = "7.7 °C"
$ws.Range("E14").Value2 = "2026-02-05 22:18:16"
$ws.Range("H14").Formula = "'74%"
$ws.Range("O14").Value2 = "-2.4 °C"
$ws.Range("E15").Value2 = "2026-02-05 22:18:19"
$ws.Range("O15").Value2 = "9.0 °C"
$ws.Range("E16").Value2 = "2026-02-05 22:18:21"
$ws.Range("O16").Value2 = "4.0 °C"
$ws.Range("E17").Value2 = "2026-02-05 22:18:24"
$ws.Range("I17").Value2 = "8.8 mm"
$ws.Range("O17").Value2 = "1.1 °C"
$ws.Range("E18").Value2 = "2026-02-05 22:18:26"
$ws.Range("E19").Value2 = "2026-02-05 22:18:29"
$ws.Range("J19").Value2 = "992.6 hPa"
$ws.Range("E20").Value2 = "2026-02-05 22:18:32"
$ws.Range("E21").Value2 = "2026-02-05 22:18:34"
$ws.Range("E22").Value2 = "2026-02-05 22:18:37"
$ws.Range("O22").Value2 = "9.4 °C"
$ws.Range("E23").Value2 = "2026-02-05 22:18:39"
$ws.Range("E24").Value2 = "2026-02-05 22:18:41"
$ws.Range("H24").Formula = "'75%"
$ws.Range("O24").Value2 = "10.6 °C"
$ws.Range("E25").Value2 = "2026-02-05 22:18:44"
$ws.Range("E26").Value2 = "2026-02-05 22:18:46"
$ws.Range("O26").Value2 = "-0.6 °C"
$ws.Range("E27").Value2 = "2026-02-05 22:18:49"
$ws.Range("E28").Value2 = "2026-02-05 22:18:51"
$ws.Range("H28").Formula = "'93%"
$ws.Range("O28").Value2 = "3.0 °C"
$ws.Range("E29").Value2 = "2026-02-05 22:18:54"
$ws.Range("H29").Formula = "'77%"
$ws.Range("O29").Value2 = "9.7 °C"
$ws.Range("E30").Value2 = "2026-02-05 22:18:56"
$ws.Range("H30").Formula = "'65%"
$ws.Range("E31").Value2 = "2026-02-05 22:18:59"
$ws.Range("I31").Value2 = "19.7 mm"
$ws.Range("E32").Value2 = "2026-02-05 22:19:01"
$ws.Range("J32").Value2 = "991.9 hPa"
$ws.Range("E33").Value2 = "2026-02-05 22:19:04"
$ws.Range("O33").Value2 = "9.7 °C"
$ws.Range("E34").Value2 = "2026-02-05 22:19:06"
$ws.Range("O34").Value2 = "4.6 °C"
$ws.Range("E35").Value2 = "2026-02-05 22:19:09"
$ws.Range("E36").Value2 = "2026-02-05 22:19:11"
